$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-12 19:48:29'
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '63%'
$ws.Range("H2").NumberFormat = "General"
$ws.Range("E3").Value = '2026-02-12 19:48:31'
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '73%'
$ws.Range("H3").NumberFormat = "General"
$ws.Range("E4").Value = '2026-02-12 19:48:33'
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = '37%'
$ws.Range("H4").NumberFormat = "General"
$ws.Range("J4").Value = '999.0 hPa'
$ws.Range("E5").Value = '2026-02-12 19:48:36'
$ws.Range("E6").Value = '2026-02-12 19:48:38'
$ws.Range("J6").Value = '998.8 hPa'
$ws.Range("E7").Value = '2026-02-12 19:48:41'
$ws.Range("J7").Value = '1001.6 hPa'
$ws.Range("K7").Value = '14.0 MJ/m2'
$ws.Range("E8").Value = '2026-02-12 19:48:44'
$ws.Range("J8").Value = '1001.0 hPa'
$ws.Range("E9").Value = '2026-02-12 19:48:46'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '68%'
$ws.Range("H9").NumberFormat = "General"
$ws.Range("E10").Value = '2026-02-12 19:48:49'
$ws.Range("O10").Value = '15.0 °C'
$ws.Range("E11").Value = '2026-02-12 19:48:51'
$ws.Range("E12").Value = '2026-02-12 19:48:54'
$ws.Range("E13").Value = '2026-02-12 19:48:56'
$ws.Range("J13").Value = '1001.5 hPa'
$ws.Range("E14").Value = '2026-02-12 19:48:58'
$ws.Range("E15").Value = '2026-02-12 19:49:01'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '54%'
$ws.Range("H15").NumberFormat = "General"
$ws.Range("E16").Value = '2026-02-12 19:49:04'
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '64%'
$ws.Range("H16").NumberFormat = "General"
$ws.Range("E17").Value = '2026-02-12 19:49:06'
$ws.Range("E18").Value = '2026-02-12 19:49:08'
$ws.Range("J18").Value = '999.2 hPa'
$ws.Range("O18").Value = '16.9 °C'
$ws.Range("E19").Value = '2026-02-12 19:49:11'
$ws.Range("O19").Value = '8.2 °C'
$ws.Range("E20").Value = '2026-02-12 19:49:13'
$ws.Range("O20").Value = '-3.4 °C'
$ws.Range("E21").Value = '2026-02-12 19:49:16'
$ws.Range("J21").Value = '1002.0 hPa'
$ws.Range("O21").Value = '9.4 °C'
$ws.Range("E22").Value = '2026-02-12 19:49:19'
$ws.Range("E23").Value = '2026-02-12 19:49:21'
$ws.Range("E24").Value = '2026-02-12 19:49:23'
$ws.Range("J24").Value = '1006.4 hPa'
$ws.Range("E25").Value = '2026-02-12 19:49:26'
$ws.Range("E26").Value = '2026-02-12 19:49:28'
$ws.Range("J26").Value = '998.4 hPa'
$ws.Range("O26").Value = '6.0 °C'
$ws.Range("E27").Value = '2026-02-12 19:49:31'
$ws.Range("O27").Value = '-1.5 °C'
$ws.Range("E28").Value = '2026-02-12 19:49:34'
$ws.Range("J28").Value = '998.6 hPa'
$ws.Range("O28").Value = '14.3 °C'
$ws.Range("E29").Value = '2026-02-12 19:49:36'
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = '56%'
$ws.Range("H29").NumberFormat = "General"
$ws.Range("N29").Value = '8.8 °C 19:27 TU'
$ws.Range("O29").Value = '14.9 °C'
$ws.Range("E30").Value = '2026-02-12 19:49:39'
$ws.Range("J30").Value = '999.1 hPa'
$ws.Range("O30").Value = '12.4 °C'
$ws.Range("E31").Value = '2026-02-12 19:49:41'
$ws.Range("J31").Value = '998.5 hPa'
$ws.Range("E32").Value = '2026-02-12 19:49:44'
$ws.Range("E33").Value = '2026-02-12 19:49:46'
$ws.Range("J33").Value = '1001.2 hPa'
$ws.Range("E34").Value = '2026-02-12 19:49:49'
$ws.Range("E35").Value = '2026-02-12 19:49:51'
$ws.Range("J35").Value = '1007.9 hPa'
$ws.Range("O35").Value = '7.9 °C'
$ws.Range("E36").Value = '2026-02-12 19:49:54'
$ws.Range("J36").Value = '999.4 hPa'
$ws.Range("K36").Value = '13.6 MJ/m2'
$ws.Range("E37").Value = '2026-02-12 19:49:56'
$ws.Range("J37").Value = '999.9 hPa'
$ws.Range("O37").Value = '10.2 °C'
$ws.Range("E38").Value = '2026-02-12 19:49:59'
$ws.Range("O38").Value = '15.9 °C'
$ws.Range("E39").Value = '2026-02-12 19:50:01'
$ws.Range("E40").Value = '2026-02-12 19:50:04'
$ws.Range("J40").Value = '1002.8 hPa'
$ws.Range("O40").Value = '9.8 °C'
$ws.Range("E41").Value = '2026-02-12 19:50:06'
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = '34%'
$ws.Range("H41").NumberFormat = "General"
$ws.Range("J41").Value = '1005.5 hPa'
$ws.Range("E42").Value = '2026-02-12 19:50:09'
$ws.Range("O42").Value = '14.4 °C'
$ws.Range("E43").Value = '2026-02-12 19:50:11'
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = '52%'
$ws.Range("H43").NumberFormat = "General"
$ws.Range("E44").Value = '2026-02-12 19:50:14'
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = '69%'
$ws.Range("H44").NumberFormat = "General"
$ws.Range("E45").Value = '2026-02-12 19:50:16'
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = '52%'
$ws.Range("H45").NumberFormat = "General"
$ws.Range("J45").Value = '1004.7 hPa'
$ws.Range("N45").Value = '3.2 °C 19:12 TU'
$ws.Range("O45").Value = '7.2 °C'
$ws.Range("E46").Value = '2026-02-12 19:50:19'
$ws.Range("J46").Value = '1007.2 hPa'
$ws.Range("N46").Value = '11.4 °C 19:13 TU'
$ws.Range("O46").Value = '16.0 °C'
